# Insert a new data row at row 43 (pushes existing rows 43-146 down to 44-147)
# and populate it with a new Murcott / Tercera price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("43:43").Insert()

$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value = "Arica y Parinacota"
$ws.Range("D43").Value = 44980
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100102
$ws.Range("H43").Value = "Cítricos"
$ws.Range("I43").Value = 100102004
$ws.Range("J43").Value = "Mandarina"
$ws.Range("K43").Value = "Murcott"
$ws.Range("L43").Value = "Tercera"
$ws.Range("M43").Value = 250
$ws.Range("N43").Value = 16000
$ws.Range("O43").Value = 17000
$ws.Range("P43").Value = 16600
$ws.Range("Q43").Value = "`$/caja 20 kilos"
$ws.Range("R43").Value = "Región de O'Higgins"
$ws.Range("S43").Value = 830
$ws.Range("T43").Value = 20
